# Cost.xlsx: "unify the conception of DataNode, DataTable, Entity."
#
# The only content-level changes the Excel object model lets us reproduce
# from the target diff are:
#   1. The worksheet is renamed from "Property1" to "DataNode".
#   2. The user's active selection on that sheet moves from A9 to D39
#      (cursor/selection position at save time).
#
# (Everything else in the upstream diff - fileVersion/rupBuild numbers,
# the absPath of the author's machine, bookViews window geometry, the
# xr/xr2/xr3/xr16/x16r2 revision-tracking namespaces and uid= GUIDs, the
# new phonetic-guide font + <phoneticPr>, the x15:timelineStyles extLst
# entry, the "Normal" -> "常规" cell style rename, and the sub-pixel
# column-width drift - are artifacts of the file being resaved by a
# different Excel build/platform and are not reachable through the
# Excel COM object model.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet.
$ws.Name = "DataNode"

# 2. Move the selection to D39 (matches the <selection .../> in the diff).
$ws.Range("D39").Select()
